$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 header values
$ws.Range("B1").Value = 15
$ws.Range("C1").Value = 16
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 16

# Row 2 values
$ws.Range("B2").Value = 16.150000000000002
$ws.Range("C2").Value = 13.750000000000002
$ws.Range("D2").Value = 13.45
$ws.Range("E2").ClearContents()

# Row 3 values
$ws.Range("B3").ClearContents()
$ws.Range("C3").Value = 6.8500000000000005
$ws.Range("D3").Value = 11.299999999999999
$ws.Range("E3").Value = 13.05

# Selection range update (matches diff's sheetView selection change)
$ws.Range("B1:E3").Select()
